# Update odds values on Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 changes
$ws.Range("G4").Value = 1.95
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 2.75
$ws.Range("L4").Value = 4.75
$ws.Range("X4").Value = 8
$ws.Range("AE4").Value = 19
$ws.Range("AI4").Value = 21
$ws.Range("AN4").Value = 3.75
$ws.Range("AO4").Value = 11
$ws.Range("AW4").Value = 6
$ws.Range("AY4").Value = 41

# Row 5 changes
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 11
